$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("N2").Value = 8.5
$ws.Range("X2").Value = 8
$ws.Range("AG2").Value = 401
$ws.Range("AL2").Value = 41
$ws.Range("AO2").Value = 10
$ws.Range("BA2").Value = 126

# Row 3 updates
$ws.Range("G3").Value = 9
$ws.Range("J3").Value = 9
$ws.Range("N3").Value = 13
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 2
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
